# Final changes of 20th May 2022
$wb = $excel.ActiveWorkbook

$wsCreation = $wb.Worksheets.Item("RTECreation")

# Update the job creation reference values on the RTECreation sheet
# (set column B before column A so new shared-string entries are appended
# in the same order the source workbook uses)
$wsCreation.Range("B2").Value = "RT00001880"
$wsCreation.Range("A2").Value = "RTE Job Creation_PreProd"
$wsCreation.Range("A3").Value = "RTE Job Creation ForCrud_PreProd"
$wsCreation.Range("B3").Value = "RT00001881"

# Make RTECreation the active sheet / active tab, with B9 selected
# (this also clears tabSelected on whichever sheet was active before,
# i.e. ShipmentDetails)
$wsCreation.Activate()
$wsCreation.Range("B9").Select()

$wb.Save()
